# Logbook sections 14 and 15.
# Applies the edits captured in the source diff:
#  - updates three existing date cells (C14:C16)
#  - fills in logbook entry #13 (row 17) - note the date was typed as plain
#    text "6/052021" (not a real date) by the original author
#  - fills in logbook entries #14 and #15 (rows 18-19) with real dates
#  - adjusts row heights for the newly-populated rows
#  - updates the sheet view (scroll position / active selection)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix up existing dates for entries #10-#12 -----------------------------
$ws.Range("C14").Value = 44257
$ws.Range("C15").Value = 44258
$ws.Range("C16").Value = 44259

# --- Entry #13 (row 17) -----------------------------------------------------
$ws.Range("B17").Value = 13
$ws.Range("C17").Value = "6/052021"
$ws.Range("D17").Value = "Began designing background image on Pisekel"
$ws.Range("E17").Value = "Background image design started on piskel. Planning to be a 700x300ish size, and take after a pixel neon city design."
$ws.Rows.Item(17).RowHeight = 29

# --- Entry #14 (row 18) -----------------------------------------------------
$ws.Range("B18").Value = 14
$ws.Range("C18").Value = 44323
$ws.Range("D18").Value = "Began testing multiprocessing"
$ws.Range("E18").Value = "Came up with the idea of using multiprocessing as a way of creating the required path images without lagging the existing game. This makes use of parralel computing by using two cores instead of one, and limits the lagging which would normally be present if the game were to generate images as part of the normal script in one process."
$ws.Rows.Item(18).RowHeight = 57.5

# --- Entry #15 (row 19) -----------------------------------------------------
$ws.Range("B19").Value = 15
$ws.Range("C19").Value = 44326
$ws.Range("D19").Value = "Began implementing image generation for path generation."
$ws.Range("E19").Value = "Created base image to be rotated to a specified angle and cropped as per the automatic path generation process. The process also now runs parralel to the existing program, instead of inline and joining the current execution thread."
$ws.Rows.Item(19).RowHeight = 43

# --- Update the sheet view: scrolled down a bit, new active selection ------
try {
    $excel.ActiveWindow.ScrollRow = 10
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
}
$ws.Range("B20").Select() | Out-Null
